$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q4" worksheet right before the existing "2022-Q3"
#    worksheet (which is currently the 2nd sheet).
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item(2)
$q4Sheet = $wb.Worksheets.Add($q3Sheet)
$q4Sheet.Name = "2022-Q4"

# Re-fetch the Q3 sheet by name since the reference above is no longer
# reliable after Add() repositions things.
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# Copy the header row (B1:H1) style from the "2022-Q3" sheet so the new
# sheet gets identical look & feel (bold, centered, bordered header cells).
$q3Sheet.Range("B1:H1").Copy($q4Sheet.Range("B1:H1"))

# The "2022-Q3" sheet only has 2 data rows, but our new sheet needs 4, so
# copy the styled index column (A2:A5, all style "2") from the "2021-Q3"
# sheet instead, which already has 4 data rows.
$q3Sheet4Rows = $wb.Worksheets.Item("2021-Q3")
$q3Sheet4Rows.Range("A2:A5").Copy($q4Sheet.Range("A2:A5"))

# Fill in the A column (row index) values.
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("A4").Value = 2
$q4Sheet.Range("A5").Value = 3

# Columns B through G hold text values (fund codes / formatted numbers),
# force text formatting so Excel does not auto-convert them to numbers.
$q4Sheet.Range("B2:G5").NumberFormat = "@"

$q4Sheet.Range("B2").Value = "014320"
$q4Sheet.Range("C2").Value = "德邦半导体产业混合C"
$q4Sheet.Range("D2").Value = "1.52"
$q4Sheet.Range("E2").Value = "92.57"
$q4Sheet.Range("F2").Value = "6.38"
$q4Sheet.Range("G2").Value = "0.0970"
$q4Sheet.Range("H2").Value = 5

$q4Sheet.Range("B3").Value = "014319"
$q4Sheet.Range("C3").Value = "德邦半导体产业混合A"
$q4Sheet.Range("D3").Value = "0.37"
$q4Sheet.Range("E3").Value = "92.57"
$q4Sheet.Range("F3").Value = "6.38"
$q4Sheet.Range("G3").Value = "0.0236"
$q4Sheet.Range("H3").Value = 5

$q4Sheet.Range("B4").Value = "016238"
$q4Sheet.Range("C4").Value = "华夏数字经济龙头混合C"
$q4Sheet.Range("D4").Value = "0.33"
$q4Sheet.Range("E4").Value = "91.36"
$q4Sheet.Range("F4").Value = "4.03"
$q4Sheet.Range("G4").Value = "0.0133"
$q4Sheet.Range("H4").Value = 8

$q4Sheet.Range("B5").Value = "016237"
$q4Sheet.Range("C5").Value = "华夏数字经济龙头混合A"
$q4Sheet.Range("D5").Value = "0.20"
$q4Sheet.Range("E5").Value = "91.36"
$q4Sheet.Range("F5").Value = "4.03"
$q4Sheet.Range("G5").Value = "0.0081"
$q4Sheet.Range("H5").Value = 8

# ---------------------------------------------------------------------------
# 2. Insert a new row 2 in the "总计" summary sheet for the 2022-Q4 figures,
#    pushing all existing rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Re-apply the index-column style (copied from the row that is now row 3)
# and clear the stray formatting Insert() leaves behind on B2:D2.
$summary.Range("A3").Copy($summary.Range("A2"))
$summary.Range("B2:D2").ClearFormats()

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.14

# Column A is just a plain 0-based row index independent of quarter content;
# re-stamp it for every data row (2022-Q4 down to 2021-Q1) after the insert.
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
